$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and
# expand the range to cover the whole paragraph.
$startRng = $d.Content
[void]$startRng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$startRng.Expand(4)

# The blank paragraph immediately preceding it (the separator paragraph
# right after the bibliography text) is removed together with it, so pull
# the deletion start back to include it.
$precedingBlank = $d.Range($startRng.Start - 1, $startRng.Start - 1)
[void]$precedingBlank.Expand(4)

# Locate the copyright/footer paragraph and expand to the whole paragraph.
$endRng = $d.Content
[void]$endRng.Find.Execute("Contact: luizeleno@usp.br", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$endRng.Expand(4)

# Delete everything from the blank separator paragraph through the end of
# the copyright paragraph (inclusive), leaving the remaining trailing blank
# paragraph and the page-break paragraph intact.
$delRange = $d.Range($precedingBlank.Start, $endRng.End)
$delRange.Delete()
